$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1123
$ws1.Range("F5").Value = 2758
$ws1.Range("F6").Value = 229
$ws1.Range("F7").Value = 692
$ws1.Range("F9").Value = 268
$ws1.Range("F11").Value = 698
$ws1.Range("F12").Value = 103
$ws1.Range("F14").Value = 1616
$ws1.Range("F17").Value = 200

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 120
$ws2.Range("C10").Value = "广州·HAG·CHINA TOUR 2024 - Phase 2 -1st Lve in Guangzhou「 声 」"
$ws2.Range("F10").Value = 24
$ws2.Range("F12").Value = 49
$ws2.Range("F18").Value = 34

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6354
$ws3.Range("F4").Value = 2020
$ws3.Range("F5").Value = 257

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6354
$ws4.Range("F4").Value = 2020
$ws4.Range("F5").Value = 257
$ws4.Range("F12").Value = 1123
$ws4.Range("F16").Value = 2758
$ws4.Range("F17").Value = 120
$ws4.Range("F18").Value = 229
$ws4.Range("C19").Value = "广州·HAG·CHINA TOUR 2024 - Phase 2 -1st Lve in Guangzhou「 声 」"
$ws4.Range("F19").Value = 24
$ws4.Range("F21").Value = 49
$ws4.Range("F22").Value = 692
$ws4.Range("F24").Value = 268
$ws4.Range("F27").Value = 698
$ws4.Range("F28").Value = 103
$ws4.Range("F31").Value = 1616
$ws4.Range("F36").Value = 200
$ws4.Range("F38").Value = 34
